$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the per-line-item columns (Item ID, Item Name, Item Type,
# Item Unit Of measure, Item MPN, Item Period, Quantity) - columns F:L -
# from the "not per line item" Fulfillment Requests report.
$ws.Range("F1:L1").EntireColumn.Delete()

# Re-apply the autofilter over the new, narrower header range so the
# worksheet's autoFilter/_FilterDatabase reference matches the shrunk data.
$ws.AutoFilterMode = $false
$ws.Range("A1:Z1").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Data!`$A`$1:`$Z`$1"
    }
}

# Restore the cursor position the author left the sheet in after editing.
$ws.Range("F3").Select()
